$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1): shift labels right by one column, with "max" wrapping to E1
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# Update data row (row 2): C2 becomes text "f__Turicibacteraceae", D2 stays the same text,
# E2 becomes numeric 1
$ws.Range("C2").Value = "f__Turicibacteraceae"
$ws.Range("D2").Value = "f__Turicibacteraceae"
$ws.Range("E2").Value = 1
